$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 6 (hunk 0)
$ws.Range("H6").Value = 185.18182
$ws.Range("I6").Value = 193.8
$ws.Range("K6").Value = 581.4000000000001
$ws.Range("M6").Value = -469.4000000000001

# row 17 (hunk 1)
$ws.Range("H17").Value = 1437.5
$ws.Range("J17").Value = 1437.5
$ws.Range("L17").Value = 4312.5
$ws.Range("N17").Value = -4648.5

# row 62 (hunk 2)
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

# row 65 (hunk 3)
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

# row 103 (hunk 4)
$ws.Range("H103").Value = 1233.8334
$ws.Range("I103").Value = 849
$ws.Range("J103").Value = 1426.25
$ws.Range("K103").Value = 2547
$ws.Range("L103").Value = 4278.75
$ws.Range("M103").Value = -1961
$ws.Range("N103").Value = -5450.75

# row 107 (hunk 5)
$ws.Range("H107").Value = 1219
$ws.Range("J107").Value = 1475.3334
$ws.Range("L107").Value = 1475.3334
$ws.Range("N107").Value = -5315.3334

# row 113 (hunk 6)
$ws.Range("H113").Value = 2499
$ws.Range("I113").Value = 998
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 998
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = 2256
$ws.Range("N113").Value = -10508

# row 125 (hunk 7)
$ws.Range("H125").Value = 2747.818
$ws.Range("J125").Value = 2244.25
$ws.Range("L125").Value = 20198.25
$ws.Range("N125").Value = -25118.25

# row 137 (hunk 8)
$ws.Range("H137").Value = 671.8
$ws.Range("I137").Value = 417.25
$ws.Range("J137").Value = 841.5
$ws.Range("K137").Value = 1251.75
$ws.Range("L137").Value = 2524.5
$ws.Range("M137").Value = 1298.25
$ws.Range("N137").Value = -7624.5

$ws = $wb.Worksheets.Item("ARM")
# row 2 (hunk 9)
$ws.Range("H2").Value = 1330.5385
$ws.Range("J2").Value = 2499.5
$ws.Range("L2").Value = 2499.5
$ws.Range("N2").Value = -2725.5

# row 110 (hunk 10)
$ws.Range("H110").Value = 1735.9333
$ws.Range("I110").Value = 1826.0769
$ws.Range("K110").Value = 1826.0769
$ws.Range("M110").Value = 218.9231

# row 116 (hunk 11)
$ws.Range("H116").Value = 1330.5385
$ws.Range("J116").Value = 2499.5
$ws.Range("L116").Value = 2499.5
$ws.Range("N116").Value = -7087.5

$ws = $wb.Worksheets.Item("BSM")
# row 3 (hunk 12)
$ws.Range("H3").Value = 1330.5385
$ws.Range("J3").Value = 2499.5
$ws.Range("L3").Value = 2499.5
$ws.Range("N3").Value = -2727.5

# row 103 (hunk 13)
$ws.Range("H103").Value = 22194.691
$ws.Range("J103").Value = 22194.691
$ws.Range("L103").Value = 22194.691
$ws.Range("N103").Value = -24538.691

# row 107 (hunk 14)
$ws.Range("H107").Value = 5519.1665
$ws.Range("I107").Value = 4823
$ws.Range("K107").Value = 4823
$ws.Range("M107").Value = -2903

# row 134 (hunk 15)
$ws.Range("H134").Value = 3328.1667
$ws.Range("I134").Value = 3406.2942
$ws.Range("K134").Value = 10218.8826
$ws.Range("M134").Value = -7683.882599999999

$ws = $wb.Worksheets.Item("CRP")
# row 58 (hunk 16)
$ws.Range("H58").Value = 4251.857
$ws.Range("I58").Value = 2841.8
$ws.Range("K58").Value = 2841.8
$ws.Range("M58").Value = -2638.8

# row 107 (hunk 17)
$ws.Range("H107").Value = 1359.6666
$ws.Range("I107").Value = 1088.7142
$ws.Range("K107").Value = 1088.7142
$ws.Range("M107").Value = 831.2858000000001

# row 122 (hunk 18)
$ws.Range("H122").Value = 2376
$ws.Range("I122").Value = 1842
$ws.Range("J122").Value = 2999
$ws.Range("K122").Value = 5526
$ws.Range("L122").Value = 8997
$ws.Range("M122").Value = -3076
$ws.Range("N122").Value = -13897

# row 134 (hunk 19)
$ws.Range("H134").Value = 3000
$ws.Range("I134").Value = 3000
$ws.Range("K134").Value = 9000
$ws.Range("M134").Value = -6465

# row 136 (hunk 20)
$ws.Range("H136").Value = 4251.857
$ws.Range("I136").Value = 2841.8
$ws.Range("K136").Value = 8525.400000000001
$ws.Range("M136").Value = -5975.400000000001

$ws = $wb.Worksheets.Item("CUL")
# row 55 (hunk 21)
$ws.Range("H55").Value = 300.5
$ws.Range("J55").Value = 300
$ws.Range("L55").Value = 900
$ws.Range("N55").Value = -1254

# row 56 (hunk 22)
$ws.Range("H56").Value = 10619.952
$ws.Range("I56").Value = 10619.952
$ws.Range("K56").Value = 10619.952
$ws.Range("M56").Value = -10089.952

# row 131 (hunk 23)
$ws.Range("H131").Value = 2495.5
$ws.Range("I131").Value = 1475.3334
$ws.Range("K131").Value = 4426.0002
$ws.Range("M131").Value = 613.9997999999996

$ws = $wb.Worksheets.Item("GSM")
# row 2 (hunk 24)
$ws.Range("H2").Value = 85.42104999999999
$ws.Range("I2").Value = 81.84614999999999
$ws.Range("K2").Value = 81.84614999999999
$ws.Range("M2").Value = 31.15385000000001

$ws = $wb.Worksheets.Item("LTW")
# row 7 (hunk 25)
$ws.Range("H7").Value = 7087
$ws.Range("I7").Value = 5000
$ws.Range("J7").Value = 7782.6665
$ws.Range("K7").Value = 5000
$ws.Range("L7").Value = 7782.6665
$ws.Range("M7").Value = -4888
$ws.Range("N7").Value = -8006.6665

# row 46 (hunk 26)
$ws.Range("H46").Value = 1891
$ws.Range("I46").Value = 2464
$ws.Range("J46").Value = 1318
$ws.Range("K46").Value = 2464
$ws.Range("L46").Value = 1318
$ws.Range("M46").Value = -2276
$ws.Range("N46").Value = -1694

# row 55 (hunk 27)
$ws.Range("H55").Value = 1927.4
$ws.Range("J55").Value = 2284
$ws.Range("L55").Value = 2284
$ws.Range("N55").Value = -2630

# row 59 (hunk 28)
$ws.Range("H59").Value = 20000
$ws.Range("J59").Value = 20000
$ws.Range("L59").Value = 20000
$ws.Range("N59").Value = -21308

# row 94 (hunk 29)
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()

# row 126 (hunk 30)
$ws.Range("H126").Value = 7087
$ws.Range("I126").Value = 5000
$ws.Range("J126").Value = 7782.6665
$ws.Range("K126").Value = 15000
$ws.Range("L126").Value = 23347.9995
$ws.Range("M126").Value = -12530
$ws.Range("N126").Value = -28287.9995

# row 136 (hunk 31)
$ws.Range("H136").Value = 2825.3333
$ws.Range("I136").Value = 3150
$ws.Range("J136").Value = 2663
$ws.Range("K136").Value = 9450
$ws.Range("L136").Value = 7989
$ws.Range("M136").Value = -6900
$ws.Range("N136").Value = -13089

$ws = $wb.Worksheets.Item("WVR")
# row 107 (hunk 32)
$ws.Range("H107").Value = 1818.8334
$ws.Range("I107").Value = 1756.375
$ws.Range("J107").Value = 1943.75
$ws.Range("K107").Value = 5269.125
$ws.Range("L107").Value = 5831.25
$ws.Range("M107").Value = -3349.125
$ws.Range("N107").Value = -9671.25

# row 132 (hunk 33)
$ws.Range("H132").Value = 2831.1667
$ws.Range("J132").Value = 2897.5
$ws.Range("L132").Value = 8692.5
$ws.Range("N132").Value = -13752.5

# row 136 (hunk 34)
$ws.Range("H136").Value = 3813.611
$ws.Range("I136").Value = 3921.3333
$ws.Range("K136").Value = 11763.9999
$ws.Range("M136").Value = -9213.999899999999
